$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The rnaSampleNumber column (F) for rows 28-53 was re-numbered to continue
# the sequence started in rows 2-27 (1..26) instead of restarting at 1, i.e.
# each value gets shifted up by 26 (1->27, 2->28, ... 26->52).
for ($row = 28; $row -le 53; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $cell.Value() + 26
}

# Update the sheet view: scroll the window down and move the active
# selection from the old blank-column selection to cell J44.
$ws.Range("J44").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
